# study ensembleModel post competition
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: ensembleModel results, highlighted with a yellow fill
$ws.Range("A6").Value = "ensembleModel"
$ws.Range("B6").Value = 0.8956273
$ws.Range("C6").Value = 0.93165
$ws.Range("D6").Value = 0.90623
$ws.Range("F6").Value = "good but spent too much time"
$ws.Rows.Item(6).Interior.Color = 65535

# Row 7: ensembleModel_gbm results
$ws.Range("A7").Value = "ensembleModel_gbm"
$ws.Range("B7").Value = 0.9001117
$ws.Range("C7").Value = 0.92801
$ws.Range("D7").Value = 0.90617
$ws.Rows.Item(7).Interior.ColorIndex = -4142

# Restore the active cell selection as it was left by the author
$ws.Range("C27").Select
